$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
# These rows previously held mixed data for several other tickers; this patch
# brings them in line with the "NOW" (ServiceNow) series already present in row 5,
# fixing the contamination caused by extra files in the data preparation step.
$data = @"
2|78.55000305175781|74.86000061035156|83.51999664306641|70.31999969482422|207519000
3|75.41999816894531|80.5|81.20999908447266|72.08000183105469|207519000
4|69.13999938964844|81.65000152587891|82.91000366210938|67.65000152587891|207519000
6|61.13000106811523|71.48000335693359|76.83999633789062|60.04999923706055|207519000
7|66.51999664306641|74.91999816894531|78.76999664306641|64.30999755859375|207519000
8|79.06999969482422|87.91000366210938|89.79000091552734|73.01999664306641|207519000
9|75.30999755859375|90.62000274658205|92.98000335693359|74.62999725341797|207519000
10|87.59999847412109|94.48000335693359|96.91999816894533|84.02999877929688|207519000
11|106.5500030517578|110.4499969482422|115.8499984741211|103.2200012207031|207519000
12|117.6500015258789|126.370002746582|129.5599975585938|116.2200012207031|207519000
13|131|148.8699951171875|150.8500061035156|130.1399993896484|207519000
14|164.4400024414062|166.1399993896484|175.5399932861328|156|207519000
15|171.6100006103516|175.9600067138672|194.8099975585937|170.3899993896484|207519000
16|197.25|181.0399932861328|199.229995727539|157.6000061035156|207519000
17|174.1199951171875|220.0200042724609|228.4100036621093|168.4700012207031|207519000
18|249.9900054931641|271.510009765625|273.8200073242188|230.1399993896484|207519000
19|281.989990234375|277.3900146484375|303.1700134277344|272.6099853515625|207519000
20|255.3200073242188|247.259994506836|276.5|213.9900054931641|207519000
21|284.9599914550781|338.2300109863281|343.8900146484375|284.2200012207031|207519000
22|274.3900146484375|351.5400085449219|357.7799987792969|246.6699981689453|207519000
23|405.3800048828125|439.2000122070313|454.7000122070313|390.8399963378906|207519000
24|491.5700073242188|497.5700073242188|533.3699951171875|478.989990234375|207519000
25|556.9400024414062|543.1599731445312|563.8800048828125|507.0599975585938|207519000
26|508|506.3699951171875|565.2899780273438|500.2799987792969|207519000
27|545|587.8900146484375|608.780029296875|539.75|207519000
28|623.780029296875|697.760009765625|698.4299926757812|609.469970703125|207519000
29|648.1599731445312|585.780029296875|651.4099731445312|479.3500061035156|207519000
30|556.0999755859375|478.1000061035156|575.8099975585938|453.9200134277344|207519000
31|474.9299926757813|446.6600036621094|506.4500122070313|407.5|207519000
32|380.0799865722656|420.739990234375|433.3900146484375|337|207519000
33|395.0799865722656|455.1300048828125|477.7300109863281|353.6199951171875|207519000
34|461.5299987792969|459.4200134277344|485.5799865722656|427.6799926757813|207519000
35|560.8499755859375|583|614.3599853515625|548.4299926757812|207519000
36|554.75|581.8499755859375|584.8400268554688|527.239990234375|207519000
37|697.0499877929688|765.4000244140625|789.9199829101562|671.219970703125|207519000
38|763.989990234375|693.3300170898438|795.4400024414062|689|207519000
39|784|814.3900146484375|850.3300170898438|729.22998046875|207519000
40|898.5900268554688|932.989990234375|979.780029296875|863.239990234375|207519000
41|1073.52001953125|1018.380004882812|1198.089965820312|990|207519000
42|798.219970703125|955.010009765625|961.030029296875|678.6599731445312|207519000
43|1027.260009765625|943.1199951171876|1057.390014648438|934.1699829101562|207519000
"@

$lines = $data -split "`r?`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $openPrice = [double]$parts[1]
    $closePrice = [double]$parts[2]
    $highPrice = [double]$parts[3]
    $lowPrice = [double]$parts[4]
    $shares = [double]$parts[5]

    $ws.Cells.Item($row, 4).Value = $openPrice
    $ws.Cells.Item($row, 5).Value = $closePrice
    $ws.Cells.Item($row, 6).Value = $highPrice
    $ws.Cells.Item($row, 7).Value = $lowPrice
    $ws.Cells.Item($row, 8).Value = $shares
    $ws.Cells.Item($row, 9).Value = "NOW"
}
